$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: force a paragraph's runs to be merged into a single run by first
# setting its Text to a throwaway placeholder (so the runtime detects an
# actual change) and then setting it to the desired final text. Re-fetch the
# paragraph object each time to avoid any stale references.
# ---------------------------------------------------------------------------
function Merge-ParagraphText {
    param($textRange, [int]$paraIndex, [string]$finalText)

    $para = $textRange.Paragraphs($paraIndex)
    $para.Text = "%%TEMP%%"
    $para = $textRange.Paragraphs($paraIndex)
    $para.Text = $finalText
}

# ===========================================================================
# Slide 10 - "Other Recommended Tools"
#   - Split the run " Toolkit)" into " Toolkit" + ")"
#   - Insert a new "Twinkle (by Leigh Dodds)" bullet after the SparqlGUI one
# ===========================================================================
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange

$sparqlPara = $tr10.Paragraphs(4)
# "SparqlGUI (dotNetRDF Toolkit)" -> split off the trailing ")" into its own run
$closeParen = $sparqlPara.Characters($sparqlPara.Length - 1, 1)
$closeParen.Text = ")"

# Insert a new paragraph right after the SparqlGUI bullet
$sparqlPara = $tr10.Paragraphs(4)
$newParaBreak = [char]13
$newParaText = $newParaBreak + "Twinkle (by Leigh Dodds)"
$sparqlPara.InsertAfter($newParaText) | Out-Null

$twinklePara = $tr10.Paragraphs(5)
$twinklePara.IndentLevel = 2

# Split "Twinkle (by Leigh Dodds)" into 3 runs:
#   "Twinkle (by Leigh " / "Dodds" / ")"
$doddsRun = $twinklePara.Characters(19, 5)
$doddsRun.Text = "Dodds"
$closeParenRun = $twinklePara.Characters(24, 1)
$closeParenRun.Text = ")"

# ===========================================================================
# Slide 2 - "About Me"
#   - Merge multi-run paragraphs into single runs
# ===========================================================================
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

Merge-ParagraphText $tr2 2 "Writing up my PhD Thesis in spare time"
Merge-ParagraphText $tr2 6 "Latest Releases"
Merge-ParagraphText $tr2 7 "0.5.1 Beta (APIs)"
Merge-ParagraphText $tr2 8 "0.4.0 Beta (Toolkit)"

# ===========================================================================
# Slide 3 - "My Life as a Semantic Web Developer"
#   - Merge multi-run paragraph into a single run
# ===========================================================================
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

Merge-ParagraphText $tr3 2 "Was frustrated by limited tooling around conversion, query prototyping and editing"

# ===========================================================================
# Slide 6 - "Editing - rdfEditor"
#   - Merge multi-run paragraphs into single runs
# ===========================================================================
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange

$windowsOnlyText = "Windows only – Mono does not support the required APIs unfortunately"
Merge-ParagraphText $tr6 3 $windowsOnlyText
Merge-ParagraphText $tr6 4 "Notepad replacement for RDF editing with syntax highlighting, auto-complete and validation capabilities"

# ===========================================================================
# Slide 8 - "Management - Store Manager"
#   - Merge multi-run title paragraph into a single run
# ===========================================================================
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(1)
$tr8 = $sh8.TextFrame.TextRange

Merge-ParagraphText $tr8 1 "Management - Store Manager"
